$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Datos actualizados a 14 de Abril de 2020 a las 22:22"

$ws.Range("A21").Value = "Asturias"
$ws.Range("B21").Value = 2096
$ws.Range("C21").Value = 487
$ws.Range("D21").Value = 1443
$ws.Range("E21").Value = 166

$ws.Range("A22").Value = "Gipuzkoa/Guipuzcoa"
$ws.Range("B22").Value = 2086
$ws.Range("C22").Value = 5193
$ws.Range("D22").Value = 5174
$ws.Range("E22").Value = 136

$ws.Range("A23").Value = "Sevilla"
$ws.Range("B23").Value = 2083
$ws.Range("C23").Value = 294
$ws.Range("D23").Value = 1611
$ws.Range("E23").Value = 178

$ws.Range("B31").Value = 1520
$ws.Range("C31").Value = 513
$ws.Range("D31").Value = 898
$ws.Range("E31").Value = 109

$ws.Range("B54").Value = 102
$ws.Range("C54").Value = 24
$ws.Range("D54").Value = 76

$ws.Range("B55").Value = 97
$ws.Range("C55").Value = 27
$ws.Range("D55").Value = 66
